# Generate Report for Handback
# Adds a new handback record (a5584dfe-a772-452b-bcae-a57c8e157f46.md) as
# row 4 on the "Overview", "zh-cn" and "de-de" worksheets, expanding each
# of their tables by one row and wiring up the corresponding hyperlinks.

$wb = $excel.ActiveWorkbook

$fileName   = "a5584dfe-a772-452b-bcae-a57c8e157f46.md"
$pathName   = "e2e\" + $fileName
$ext        = ".md"
$status     = "Handed back: in sync with en-US"
$genDate    = "2016-09-07 13:07:39"

$srcCommit  = "0c933b45afb7e56003666ccc5949ceb31e1039c8"
$zhCommit   = "95b419ef0308f28ed7605bbcfc4e378534815136"
$deCommit   = "3a698e2c08cf35590c987aff6cdbe60980f08554"

$srcUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$srcCommit/e2e/$fileName"
$zhUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/$zhCommit/e2e/$fileName"
$deUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/$deCommit/e2e/$fileName"

$zhXlf       = "a5584dfe-a772-452b-bcae-a57c8e157f46.ac35884a4aa96b31ee522ff49afd0d7e0ffa5b1a.zh-cn.xlf"
$zhHandoffDt = "2016-09-07 13:07:28"
$zhHandbkDt  = "2016-09-07 13:08:28"

$deXlf       = "a5584dfe-a772-452b-bcae-a57c8e157f46.ac35884a4aa96b31ee522ff49afd0d7e0ffa5b1a.de-de.xlf"
$deHandoffDt = "2016-09-07 13:07:39"
$deHandbkDt  = "2016-09-07 13:08:47"

# ---------------------------------------------------------------------
# Overview sheet (row 4)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item("Overview")
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Cells.Item(4, 1).Value = $fileName
$wsOverview.Cells.Item(4, 3).Value = $ext
$wsOverview.Cells.Item(4, 5).Value = $status
$wsOverview.Cells.Item(4, 6).Value = $status
$wsOverview.Cells.Item(4, 7).Value = $genDate
$wsOverview.Cells.Item(4, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), $srcUrl, "", "", $pathName) | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet (row 4)
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item("zh-cn")
$loZh.ListRows.Add() | Out-Null

$wsZh.Cells.Item(4, 2).Value  = $ext
$wsZh.Cells.Item(4, 3).Value  = $status
$wsZh.Cells.Item(4, 4).Value  = "e2e"
$wsZh.Cells.Item(4, 5).Value  = "ht"
$wsZh.Cells.Item(4, 6).Value  = "True"
$wsZh.Cells.Item(4, 7).Value  = $zhXlf
$wsZh.Cells.Item(4, 8).Value  = $zhHandoffDt
$wsZh.Cells.Item(4, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(4, 10).Value = $zhXlf
$wsZh.Cells.Item(4, 11).Value = $zhHandbkDt
$wsZh.Cells.Item(4, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(4, 13).Value = "True"
$wsZh.Cells.Item(4, 15).Value = "False"

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $srcUrl, "", "", $fileName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), $zhUrl, "", "", $fileName) | Out-Null

# ---------------------------------------------------------------------
# de-de sheet (row 4)
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item("de-de")
$loDe.ListRows.Add() | Out-Null

$wsDe.Cells.Item(4, 2).Value  = $ext
$wsDe.Cells.Item(4, 3).Value  = $status
$wsDe.Cells.Item(4, 4).Value  = "e2e"
$wsDe.Cells.Item(4, 5).Value  = "ht"
$wsDe.Cells.Item(4, 6).Value  = "True"
$wsDe.Cells.Item(4, 7).Value  = $deXlf
$wsDe.Cells.Item(4, 8).Value  = $deHandoffDt
$wsDe.Cells.Item(4, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(4, 10).Value = $deXlf
$wsDe.Cells.Item(4, 11).Value = $deHandbkDt
$wsDe.Cells.Item(4, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(4, 13).Value = "True"
$wsDe.Cells.Item(4, 15).Value = "False"

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $srcUrl, "", "", $fileName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), $deUrl, "", "", $fileName) | Out-Null
